$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B2").Value = 42986527.03819948
$ws.Range("C2").Value = 59849964.84684038
$ws.Range("D2").Value = 264616988.29282954

$ws.Range("B3").Value = 50357683.38647374
$ws.Range("C3").Value = 492092848.47966754
$ws.Range("D3").Value = 14197923632.578537

$ws.Range("B4").Value = 118896424.35539998
$ws.Range("C4").Value = 300916021.42223704
$ws.Range("D4").Value = 6548654199.625336
